# Daily attendance processing - 2026-01-15 21:37:12
# Reorder the "Recorded By" (column G) email lists for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"
    3  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"
    4  = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    5  = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
    6  = "majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg"
    7  = "lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
    8  = "NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
    12 = "Madeha.Saeed@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
    13 = "esraa.mostafa@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
    19 = "Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
    27 = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
    30 = "shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
